$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2242.2173
$ws.Range("I6").Value = 1131
$ws.Range("K6").Value = 3393
$ws.Range("M6").Value = -3281
$ws.Range("H40").Value = 2291.5833
$ws.Range("I40").Value = 1999.8334
$ws.Range("K40").Value = 1999.8334
$ws.Range("M40").Value = -1824.8334
$ws.Range("H58").Value = 1571
$ws.Range("I58").Value = 656.5
$ws.Range("J58").Value = 3400
$ws.Range("K58").Value = 1969.5
$ws.Range("L58").Value = 10200
$ws.Range("M58").Value = -1819.5
$ws.Range("N58").Value = -10500
$ws.Range("H61").Value = 738.3333
$ws.Range("I61").Value = 738.3333
$ws.Range("K61").Value = 2214.9999
$ws.Range("M61").Value = -2042.9999
$ws.Range("H70").Value = 4709.375
$ws.Range("J70").Value = 5779.4165
$ws.Range("L70").Value = 17338.2495
$ws.Range("N70").Value = -17878.2495
$ws.Range("H73").Value = 4709.375
$ws.Range("J73").Value = 5779.4165
$ws.Range("L73").Value = 17338.2495
$ws.Range("N73").Value = -19210.2495
$ws.Range("H112").Value = 2281.6365
$ws.Range("J112").Value = 2281.6365
$ws.Range("L112").Value = 6844.9095
$ws.Range("N112").Value = -9060.9095
$ws.Range("H125").Value = 755
$ws.Range("I125").Value = 755
$ws.Range("K125").Value = 6795
$ws.Range("M125").Value = -4335
$ws.Range("H132").Value = 2523.5
$ws.Range("I132").Value = 2523.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7570.5
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -5040.5
$ws.Range("H135").Value = 817.4545000000001
$ws.Range("I135").Value = 630.26086
$ws.Range("K135").Value = 5672.34774
$ws.Range("M135").Value = -3137.34774
$ws.Range("H137").Value = 2088.5715
$ws.Range("I137").Value = 1551.9474
$ws.Range("J137").Value = 3221.4443
$ws.Range("K137").Value = 4655.8422
$ws.Range("L137").Value = 9664.332900000001
$ws.Range("M137").Value = -2105.8422
$ws.Range("N137").Value = -14764.3329

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4570.7017
$ws.Range("I32").Value = 2810
$ws.Range("K32").Value = 2810
$ws.Range("M32").Value = -2523
$ws.Range("H74").Value = 1320.6389
$ws.Range("I74").Value = 878.3125
$ws.Range("J74").Value = 4859.25
$ws.Range("K74").Value = 878.3125
$ws.Range("L74").Value = 4859.25
$ws.Range("M74").Value = -4.3125
$ws.Range("N74").Value = -6607.25
$ws.Range("H77").Value = 1320.6389
$ws.Range("I77").Value = 878.3125
$ws.Range("J77").Value = 4859.25
$ws.Range("K77").Value = 4391.5625
$ws.Range("L77").Value = 24296.25
$ws.Range("M77").Value = -23.5625
$ws.Range("N77").Value = -33032.25
$ws.Range("H122").Value = 373224
$ws.Range("I122").Value = 558400.25
$ws.Range("J122").Value = 2871.5557
$ws.Range("K122").Value = 1675200.75
$ws.Range("L122").Value = 8614.667099999999
$ws.Range("M122").Value = -1672750.75
$ws.Range("N122").Value = -13514.6671
$ws.Range("H132").Value = 1672.8235
$ws.Range("I132").Value = 1672.8235
$ws.Range("K132").Value = 5018.470499999999
$ws.Range("M132").Value = -2488.470499999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 697.9091
$ws.Range("I22").Value = 697.9091
$ws.Range("K22").Value = 697.9091
$ws.Range("M22").Value = -524.9091

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3194.3333
$ws.Range("I58").Value = 1399.8572
$ws.Range("J58").Value = 5126.846
$ws.Range("K58").Value = 1399.8572
$ws.Range("L58").Value = 5126.846
$ws.Range("M58").Value = -1196.8572
$ws.Range("N58").Value = -5532.846
$ws.Range("H134").Value = 2548.875
$ws.Range("I134").Value = 2068.0386
$ws.Range("J134").Value = 3441.8572
$ws.Range("K134").Value = 6204.1158
$ws.Range("L134").Value = 10325.5716
$ws.Range("M134").Value = -3669.1158
$ws.Range("N134").Value = -15395.5716
$ws.Range("H136").Value = 3194.3333
$ws.Range("I136").Value = 1399.8572
$ws.Range("J136").Value = 5126.846
$ws.Range("K136").Value = 4199.571599999999
$ws.Range("L136").Value = 15380.538
$ws.Range("M136").Value = -1649.571599999999
$ws.Range("N136").Value = -20480.538

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 14285817
$ws.Range("I7").Value = 20000094
$ws.Range("J7").Value = 125
$ws.Range("K7").Value = 60000282
$ws.Range("L7").Value = 375
$ws.Range("M7").Value = -60000170
$ws.Range("N7").Value = -599
$ws.Range("H17").Value = 1185.1428
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = 2036.5
$ws.Range("K17").Value = 150
$ws.Range("L17").Value = 6109.5
$ws.Range("M17").Value = 19
$ws.Range("N17").Value = -6447.5
$ws.Range("H80").Value = 4099.5
$ws.Range("J80").Value = 6000
$ws.Range("L80").Value = 18000
$ws.Range("N80").Value = -19872
$ws.Range("H83").Value = 4099.5
$ws.Range("J83").Value = 6000
$ws.Range("L83").Value = 54000
$ws.Range("N83").Value = -63360
$ws.Range("H107").Value = 414.34482
$ws.Range("I107").Value = 317.16666
$ws.Range("J107").Value = 439.69565
$ws.Range("K107").Value = 951.4999799999999
$ws.Range("L107").Value = 1319.08695
$ws.Range("M107").Value = 968.5000200000001
$ws.Range("N107").Value = -5159.08695

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 964
$ws.Range("J41").Value = 978.5
$ws.Range("L41").Value = 978.5
$ws.Range("N41").Value = -1688.5
$ws.Range("H52").Value = 8750
$ws.Range("I52").Value = 5000
$ws.Range("K52").Value = 5000
$ws.Range("M52").Value = -4741
$ws.Range("H80").Value = 3998.2222
$ws.Range("I80").Value = 3305.8462
$ws.Range("J80").Value = 5798.4
$ws.Range("K80").Value = 3305.8462
$ws.Range("L80").Value = 5798.4
$ws.Range("M80").Value = -2307.8462
$ws.Range("N80").Value = -7794.4
$ws.Range("H83").Value = 3998.2222
$ws.Range("I83").Value = 3305.8462
$ws.Range("J83").Value = 5798.4
$ws.Range("K83").Value = 16529.231
$ws.Range("L83").Value = 28992
$ws.Range("M83").Value = -11537.231
$ws.Range("N83").Value = -38976
$ws.Range("H113").Value = 6225
$ws.Range("I113").Value = 6890
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 6890
$ws.Range("L113").Value = 2900
$ws.Range("M113").Value = -4720
$ws.Range("N113").Value = -7240
$ws.Range("H122").Value = 64994.188
$ws.Range("J122").Value = 168732.83
$ws.Range("L122").Value = 506198.49
$ws.Range("N122").Value = -511098.49
$ws.Range("H123").Value = 42793.145
$ws.Range("J123").Value = 42793.145
$ws.Range("L123").Value = 42793.145
$ws.Range("N123").Value = -47693.145
$ws.Range("H126").Value = 3809.8333
$ws.Range("I126").Value = 2961.25
$ws.Range("J126").Value = 5507
$ws.Range("K126").Value = 8883.75
$ws.Range("L126").Value = 16521
$ws.Range("M126").Value = -6413.75
$ws.Range("N126").Value = -21461
$ws.Range("H132").Value = 2477.3845
$ws.Range("I132").Value = 1901.875
$ws.Range("K132").Value = 5705.625
$ws.Range("M132").Value = -3175.625

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3319.5334
$ws.Range("I7").Value = 3056.6428
$ws.Range("K7").Value = 3056.6428
$ws.Range("M7").Value = -2944.6428
$ws.Range("H22").Value = 19199.666
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 28499.5
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 28499.5
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -29089.5
$ws.Range("H27").Value = 19199.666
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 28499.5
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 28499.5
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -28713.5
$ws.Range("H122").Value = 6581.2856
$ws.Range("I122").Value = 7737
$ws.Range("J122").Value = 5040.3335
$ws.Range("K122").Value = 23211
$ws.Range("L122").Value = 15121.0005
$ws.Range("M122").Value = -20761
$ws.Range("N122").Value = -20021.0005
$ws.Range("H126").Value = 3319.5334
$ws.Range("I126").Value = 3056.6428
$ws.Range("K126").Value = 9169.928400000001
$ws.Range("M126").Value = -6699.928400000001
$ws.Range("H132").Value = 3363.394
$ws.Range("I132").Value = 2966.4167
$ws.Range("K132").Value = 8899.250100000001
$ws.Range("M132").Value = -6369.250100000001
$ws.Range("H133").Value = 105000
$ws.Range("J133").Value = 105000
$ws.Range("L133").Value = 105000
$ws.Range("N133").Value = -110060
$ws.Range("H136").Value = 5398.4
$ws.Range("I136").Value = 5156
$ws.Range("K136").Value = 15468
$ws.Range("M136").Value = -12918

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10000
$ws.Range("I81").Value = 19000
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 38000
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -36939
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 10000
$ws.Range("I84").Value = 19000
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 190000
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -184696
$ws.Range("N84").Value = -20608
$ws.Range("H123").Value = 48999.5
$ws.Range("J123").Value = 48999.5
$ws.Range("L123").Value = 48999.5
$ws.Range("N123").Value = -58799.5
$ws.Range("H126").Value = 1956.8
$ws.Range("I126").Value = 1581.6666
$ws.Range("J126").Value = 5333
$ws.Range("K126").Value = 4744.9998
$ws.Range("L126").Value = 15999
$ws.Range("M126").Value = -2274.9998
$ws.Range("N126").Value = -20939
$ws.Range("H132").Value = 46017.617
$ws.Range("I132").Value = 53365.223
$ws.Range("K132").Value = 160095.669
$ws.Range("M132").Value = -157565.669
$ws.Range("H136").Value = 1185.4828
$ws.Range("I136").Value = 1195.9615
$ws.Range("K136").Value = 3587.8845
$ws.Range("M136").Value = -1037.8845
